$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "pH probe" labels, cycling through the 6-item sequence for the 12 data rows
$probeValues = @(
    "F.0.1.13_1",
    "F.0.1.13_2",
    "F.0.1.13_3",
    "F.0.1.13_4",
    "F.0.1.21_1",
    "F.0.1.21_2",
    "F.0.1.13_1",
    "F.0.1.13_2",
    "F.0.1.13_3",
    "F.0.1.13_4",
    "F.0.1.21_1",
    "F.0.1.21_2"
)

# Update existing rows 2-6 (only the "pH probe" column C changes)
for ($i = 0; $i -lt 5; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 3).Value = $probeValues[$i]
}

# Add new rows 7-13, matching the formatting/style of row 6 (the prototype data row)
for ($row = 7; $row -le 13; $row++) {
    $ws.Range("A6:H6").Copy()
    $ws.Range("A$row`:H$row").PasteSpecial(-4122)  # xlPasteFormats

    $ws.Cells.Item($row, 1).Value = $row - 1
    $ws.Cells.Item($row, 2).Value = 1
    $ws.Cells.Item($row, 3).Value = $probeValues[$row - 2]
    $ws.Cells.Item($row, 4).Value = 1440
    $ws.Cells.Item($row, 5).Value = 5.6
    $ws.Cells.Item($row, 6).Value = 6.8
    $ws.Cells.Item($row, 7).Value = 5
    $ws.Cells.Item($row, 8).Value = 1
}

$excel.CutCopyMode = 0
$null = $ws.Range("C2").Select()
